# Adds a "CODIGO" (code) column to the product table, inserting it to the
# left of "NOMBRE PRODUCTO" (rows 7-11, which previously spanned merged
# cells B:F). Column B becomes its own single-cell "CODIGO" column and the
# previous B:F merged block becomes a C:F merged block holding the product
# name, mirroring the already-existing MARCA / PRECIO / CANTIDAD columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel constant values (avoid relying on the interop enum types being
# importable by name in every host):
$xlEdgeLeft   = 7
$xlEdgeTop    = 8
$xlEdgeBottom = 9
$xlEdgeRight  = 10
$xlThin       = 2
$xlContinuous = 1
$xlNone       = -4142
$xlHAlignCenter = -4108
$xlThemeColorAccent4 = 9   # -> OOXML theme index 8, the existing header blue fill

function Set-BoxBorders {
    param($range, [bool]$left, [bool]$right)

    # Top and bottom are always present on every cell of the range.
    $range.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
    $range.Borders.Item($xlEdgeTop).Weight = $xlThin
    $range.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
    $range.Borders.Item($xlEdgeBottom).Weight = $xlThin

    if ($left) {
        $range.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
        $range.Borders.Item($xlEdgeLeft).Weight = $xlThin
    } else {
        $range.Borders.Item($xlEdgeLeft).LineStyle = $xlNone
    }

    if ($right) {
        $range.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
        $range.Borders.Item($xlEdgeRight).Weight = $xlThin
    } else {
        $range.Borders.Item($xlEdgeRight).LineStyle = $xlNone
    }
}

# Codes to attach to each product row, and the row numbers they belong to.
$codigos = @{ 8 = "SSD41"; 9 = "SSD42"; 10 = "SER74"; 11 = "SXX48" }

# ---- Header row (row 7) -----------------------------------------------
$ws.Range("B7:F7").UnMerge()

$ws.Range("C7").Value = $ws.Range("B7").Value2   # "NOMBRE PRODUCTO" -> C7
$ws.Range("B7").Value = "CÓDIGO"

$hdrCodeCell = $ws.Range("B7")
$hdrCodeCell.Interior.ThemeColor = $xlThemeColorAccent4
$hdrCodeCell.Interior.TintAndShade = 0
$hdrCodeCell.HorizontalAlignment = $xlHAlignCenter
Set-BoxBorders $hdrCodeCell $true $true

$hdrNameRange = $ws.Range("C7:F7")
$hdrNameRange.Interior.ThemeColor = $xlThemeColorAccent4
$hdrNameRange.Interior.TintAndShade = 0
$ws.Range("C7").HorizontalAlignment = $xlHAlignCenter
Set-BoxBorders $hdrNameRange $false $true

$ws.Range("C7:F7").Merge() | Out-Null

# ---- Data rows (rows 8-11) ---------------------------------------------
foreach ($r in 8..11) {
    $rowRange = "B" + $r + ":F" + $r
    $ws.Range($rowRange).UnMerge()

    $nameCell = $ws.Range("C$r")
    $codeCell = $ws.Range("B$r")

    $nameCell.Value = $codeCell.Value2   # move the product name from B to C
    $codeCell.Value = $codigos[$r]

    $codeCell.Interior.Pattern = $xlNone
    $codeCell.HorizontalAlignment = $xlHAlignCenter
    Set-BoxBorders $codeCell $true $true

    $nameRange = $ws.Range("C" + $r + ":F" + $r)
    $nameRange.Interior.Pattern = $xlNone
    $nameCell.HorizontalAlignment = $xlHAlignCenter
    Set-BoxBorders $nameRange $false $true

    $ws.Range($rowRange.Replace("B","C")) | Out-Null
    $ws.Range("C" + $r + ":F" + $r).Merge() | Out-Null
}

# ---- Misc: restore the active cell selection as recorded in the saved
# workbook (cosmetic, matches what Excel persists when a user clicks D17).
$ws.Range("D17").Select() | Out-Null
